# Generate Report for handback
# Update the "Correspond Handoff Datetime" (D) and "Correspond Handback DateTime" (G)
# columns for the d51fb795... row (row 3) on the zh-cn sheet, and for the same row
# on the de-de sheet, to reflect the newly generated handback report timestamps.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-01-08 13:33:22"
$wsZhCn.Range("G3").Value = "2016-01-08 13:34:12"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-01-08 13:33:35"
$wsDeDe.Range("G3").Value = "2016-01-08 13:34:33"
